$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 538; this shifts the existing rows 538-608
# down to 539-609 and extends the used range to A1:R609.
$ws.Rows(538).Insert()

# Populate the newly inserted row 538 with the new record's data.
$ws.Range("A538").Value = 3
$ws.Range("B538").Value = "Femacal de La Calera"
$ws.Range("C538").Value = "Coquimbo"
$ws.Range("D538").Value = 45154
$ws.Range("D538").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E538").Value = 5
$ws.Range("F538").Value = 100114013
$ws.Range("G538").Value = "Zanahoria"
$ws.Range("H538").Value = "Sin especificar"
$ws.Range("I538").Value = "Primera"
$ws.Range("J538").Value = 280
$ws.Range("K538").Value = 7000
$ws.Range("L538").Value = 7500
$ws.Range("M538").Value = 7286
$ws.Range("N538").Value = "$/saco 20 kilos"
$ws.Range("O538").Value = "Provincia de Quillota"
$ws.Range("P538").Value = 364
$ws.Range("Q538").Value = 20
$ws.Range("R538").Value = "Hortaliza"
